# Auto-applies the two-digit division answer table updates per the commit diff.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("98÷9=10, 8", $true, $false, $false, $false, $false, $true, 1, $false, "40÷2=20, 0", 2)

$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("71÷8=8, 7", $true, $false, $false, $false, $false, $true, 1, $false, "11÷5=2, 1", 2)

$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("15÷8=1, 7", $true, $false, $false, $false, $false, $true, 1, $false, "57÷4=14, 1", 2)

$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("79÷7=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "41÷4=10, 1", 2)

$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("28÷7=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "84÷6=14, 0", 2)

$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("78÷8=9, 6", $true, $false, $false, $false, $false, $true, 1, $false, "57÷8=7, 1", 2)

$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("30÷2=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "10÷8=1, 2", 2)

$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("72÷7=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "22÷8=2, 6", 2)

$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("87÷7=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "39÷5=7, 4", 2)

$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("18÷4=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "19÷8=2, 3", 2)

$cell = $t.Cell(9, 1)
$cell.Range.Find.Execute("33÷4=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "60÷9=6, 6", 2)

$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("16÷4=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "99÷8=12, 3", 2)

$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("94÷7=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "34÷9=3, 7", 2)

$cell = $t.Cell(9, 4)
$cell.Range.Find.Execute("79÷4=19, 3", $true, $false, $false, $false, $false, $true, 1, $false, "79÷9=8, 7", 2)

$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("86÷6=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "45÷8=5, 5", 2)

$cell = $t.Cell(13, 1)
$cell.Range.Find.Execute("78÷2=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "56÷5=11, 1", 2)

$cell = $t.Cell(13, 2)
$cell.Range.Find.Execute("99÷6=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "98÷9=10, 8", 2)

$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("62÷8=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "63÷5=12, 3", 2)

$cell = $t.Cell(13, 4)
$cell.Range.Find.Execute("17÷9=1, 8", $true, $false, $false, $false, $false, $true, 1, $false, "84÷6=14, 0", 2)

$cell = $t.Cell(13, 5)
$cell.Range.Find.Execute("74÷8=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "35÷2=17, 1", 2)

$cell = $t.Cell(17, 1)
$cell.Range.Find.Execute("53÷2=26, 1", $true, $false, $false, $false, $false, $true, 1, $false, "26÷3=8, 2", 2)

$cell = $t.Cell(17, 2)
$cell.Range.Find.Execute("20÷6=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "41÷3=13, 2", 2)

$cell = $t.Cell(17, 3)
$cell.Range.Find.Execute("99÷3=33, 0", $true, $false, $false, $false, $false, $true, 1, $false, "76÷8=9, 4", 2)

$cell = $t.Cell(17, 4)
$cell.Range.Find.Execute("86÷4=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "52÷8=6, 4", 2)

$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("20÷5=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "42÷7=6, 0", 2)
